# Update feed logs: append two new log rows (run_id 99 and 100) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A100").Value = 99
$ws.Range("B100").Value = 1
$ws.Range("C100").Value = "2024-06-17 00:59:34"
$ws.Range("D100").Value = 200
$ws.Range("E100").Value = 14

$ws.Range("A101").Value = 100
$ws.Range("B101").Value = 2
$ws.Range("C101").Value = "2024-06-17 00:59:35"
$ws.Range("D101").Value = 200
$ws.Range("E101").Value = 0
